# Update the PSII trends data table (Sheet1) with the new partitioned
# pesticide-concentration values that feed the PSII trends bar chart's
# (partitioned) value axis, then leave the selection positioned on the
# last updated cell (matches the scrolled/selected view captured in the
# saved workbook).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 32 (2009-2010 total row for the Burdekin / "Site A" block)
$ws.Range("D32").Value = 0
$ws.Range("E32").Value = 0
$ws.Range("F32").Value = 0
$ws.Range("G32").Value = 100
$ws.Range("H32").Value = 0
$ws.Range("I32").Value = 0
$ws.Range("J32").Value = 0

# Row 57 (2009-2010 total row for the Mackay Whitsunday / "Site C" block)
$ws.Range("D57").Value = 2
$ws.Range("E57").Value = 28
$ws.Range("F57").Value = 0
$ws.Range("G57").Value = 460
$ws.Range("H57").Value = 20
$ws.Range("I57").Value = 0
$ws.Range("J57").Value = 0

# Scroll/select so the view matches what was captured when the edit was saved.
[void]$ws.Range("A18").Select()
[void]$ws.Range("J58").Select()
